$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 7, 8, 9, 10 have their Fecha/Volumen/Precio columns rotated
# following the cycle: 2 -> 9 -> 8 -> 7 -> 10 -> 2
# (the old data in row 2 moves to row 9, old row 9 -> row 8,
#  old row 8 -> row 7, old row 7 -> row 10, old row 10 -> row 2)

$rows = @(2, 7, 8, 9, 10)
$cols = @("D", "J", "K", "L", "M", "P")

# Capture old values for each row/column combination (use Value2 to get
# raw numeric values, avoiding date-formatted strings for column D)
$old = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2()
    }
    $old[$r] = $rowVals
}

# Destination row for each source row's old data
$dest = @{ 2 = 9; 9 = 8; 8 = 7; 7 = 10; 10 = 2 }

foreach ($src in $rows) {
    $target = $dest[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value2 = $old[$src][$c]
    }
}
